$d = $word.ActiveDocument
$count = 0
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.ContextualSpacing = $false
    $count = $count + 1
}
Write-Host "Processed $count paragraphs"
